# Scheduled market-data refresh: update cached price/profit figures
# (currentAveragePrice*, LevePrice*, LeveProfit*) on the per-class Leve
# tables. Pure data overwrite - no formulas, no structural changes.
$wb = $excel.ActiveWorkbook

# --- ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(4, 8).Value = 190
$ws.Cells.Item(4, 9).Value = 190
$ws.Cells.Item(4, 11).Value = 190
$ws.Cells.Item(4, 13).Value = -76
$ws.Cells.Item(5, 8).Value = 95.666664
$ws.Cells.Item(5, 9).Value = 107.14286
$ws.Cells.Item(5, 10).Value = 55.5
$ws.Cells.Item(5, 11).Value = 107.14286
$ws.Cells.Item(5, 12).Value = 55.5
$ws.Cells.Item(5, 13).Value = 7.857140000000001
$ws.Cells.Item(5, 14).Value = -285.5
$ws.Cells.Item(19, 8).Value = 824.61536
$ws.Cells.Item(19, 10).Value = 585.5
$ws.Cells.Item(19, 12).Value = 585.5
$ws.Cells.Item(19, 14).Value = -935.5
$ws.Cells.Item(29, 8).Value = 4112
$ws.Cells.Item(29, 9).Value = 1409.6
$ws.Cells.Item(29, 11).Value = 4228.799999999999
$ws.Cells.Item(29, 13).Value = -3947.799999999999
$ws.Cells.Item(31, 8).Value = 309.33334
$ws.Cells.Item(31, 9).Value = 309.33334
$ws.Cells.Item(31, 11).Value = 928.0000200000001
$ws.Cells.Item(31, 13).Value = -698.0000200000001
$ws.Cells.Item(40, 8).Value = 4974.75
$ws.Cells.Item(40, 10).Value = 8799.429
$ws.Cells.Item(40, 12).Value = 8799.429
$ws.Cells.Item(40, 14).Value = -9149.429
$ws.Cells.Item(64, 8).Value = 8613.714
$ws.Cells.Item(64, 9).Value = 8374.25
$ws.Cells.Item(64, 10).Value = 8933
$ws.Cells.Item(64, 11).Value = 8374.25
$ws.Cells.Item(64, 12).Value = 8933
$ws.Cells.Item(64, 13).Value = -8126.25
$ws.Cells.Item(64, 14).Value = -9429
$ws.Cells.Item(67, 8).Value = 8613.714
$ws.Cells.Item(67, 9).Value = 8374.25
$ws.Cells.Item(67, 10).Value = 8933
$ws.Cells.Item(67, 11).Value = 8374.25
$ws.Cells.Item(67, 12).Value = 8933
$ws.Cells.Item(67, 13).Value = -7516.25
$ws.Cells.Item(67, 14).Value = -10649
$ws.Cells.Item(76, 8).Value = 2619.6
$ws.Cells.Item(76, 9).Value = 4999
$ws.Cells.Item(76, 10).Value = 2024.75
$ws.Cells.Item(76, 11).Value = 4999
$ws.Cells.Item(76, 12).Value = 2024.75
$ws.Cells.Item(76, 13).Value = -4684
$ws.Cells.Item(76, 14).Value = -2654.75
$ws.Cells.Item(79, 8).Value = 2619.6
$ws.Cells.Item(79, 9).Value = 4999
$ws.Cells.Item(79, 10).Value = 2024.75
$ws.Cells.Item(79, 11).Value = 4999
$ws.Cells.Item(79, 12).Value = 2024.75
$ws.Cells.Item(79, 13).Value = -3907
$ws.Cells.Item(79, 14).Value = -4208.75
$ws.Cells.Item(112, 8).Value = 1814.8235
$ws.Cells.Item(112, 9).Value = 1512.5
$ws.Cells.Item(112, 11).Value = 4537.5
$ws.Cells.Item(112, 13).Value = -3429.5

# --- ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(4, 8).Value = 1109
$ws.Cells.Item(4, 10).Value = 1109
$ws.Cells.Item(4, 12).Value = 1109
$ws.Cells.Item(4, 14).Value = -1341
$ws.Cells.Item(5, 8).Value = 72.625
$ws.Cells.Item(5, 9).Value = 63.5
$ws.Cells.Item(5, 11).Value = 63.5
$ws.Cells.Item(5, 13).Value = 48.5
$ws.Cells.Item(39, 8).Value = 3083.2
$ws.Cells.Item(39, 9).Value = 1354
$ws.Cells.Item(39, 11).Value = 1354
$ws.Cells.Item(39, 13).Value = -834
$ws.Cells.Item(61, 8).Value = 7226.091
$ws.Cells.Item(61, 9).Value = 5162.6665
$ws.Cells.Item(61, 11).Value = 5162.6665
$ws.Cells.Item(61, 13).Value = -4950.6665
$ws.Cells.Item(88, 8).Value = 1506.2
$ws.Cells.Item(88, 10).Value = 1507
$ws.Cells.Item(88, 12).Value = 1507
$ws.Cells.Item(88, 14).Value = -2319
$ws.Cells.Item(91, 8).Value = 1506.2
$ws.Cells.Item(91, 10).Value = 1507
$ws.Cells.Item(91, 12).Value = 1507
$ws.Cells.Item(91, 14).Value = -4315
$ws.Cells.Item(122, 8).Value = 923.125
$ws.Cells.Item(122, 9).Value = 798.5
$ws.Cells.Item(122, 10).Value = 1297
$ws.Cells.Item(122, 11).Value = 2395.5
$ws.Cells.Item(122, 12).Value = 3891
$ws.Cells.Item(122, 13).Value = 54.5
$ws.Cells.Item(122, 14).Value = -8791
$ws.Cells.Item(132, 8).Value = 3376.4443
$ws.Cells.Item(132, 9).Value = 3376.4443
$ws.Cells.Item(132, 11).Value = 10129.3329
$ws.Cells.Item(132, 13).Value = -7599.332900000001
$ws.Cells.Item(136, 8).Value = 7226.091
$ws.Cells.Item(136, 9).Value = 5162.6665
$ws.Cells.Item(136, 11).Value = 15487.9995
$ws.Cells.Item(136, 13).Value = -12937.9995

# --- BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(4, 8).Value = 72.625
$ws.Cells.Item(4, 9).Value = 63.5
$ws.Cells.Item(4, 11).Value = 63.5
$ws.Cells.Item(4, 13).Value = 51.5

# --- CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(7, 8).Value = 494.3125
$ws.Cells.Item(7, 9).Value = 96.583336
$ws.Cells.Item(7, 11).Value = 96.583336
$ws.Cells.Item(7, 13).Value = 16.416664
$ws.Cells.Item(22, 8).Value = 255
$ws.Cells.Item(22, 10).Value = 255
$ws.Cells.Item(22, 12).Value = 255
$ws.Cells.Item(22, 14).Value = -955
$ws.Cells.Item(31, 8).Value = 5152.4194
$ws.Cells.Item(31, 9).Value = 2393.5881
$ws.Cells.Item(31, 11).Value = 2393.5881
$ws.Cells.Item(31, 13).Value = -2098.5881
$ws.Cells.Item(34, 8).Value = 5152.4194
$ws.Cells.Item(34, 9).Value = 2393.5881
$ws.Cells.Item(34, 11).Value = 2393.5881
$ws.Cells.Item(34, 13).Value = -2191.5881

# --- CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(96, 8).Value = 3250
$ws.Cells.Item(96, 10).Value = 3250
$ws.Cells.Item(96, 12).Value = 9750
$ws.Cells.Item(96, 14).Value = -13868

# --- GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(2, 8).Value = 357.625
$ws.Cells.Item(2, 9).Value = 309.33334
$ws.Cells.Item(2, 10).Value = 386.6
$ws.Cells.Item(2, 11).Value = 309.33334
$ws.Cells.Item(2, 12).Value = 386.6
$ws.Cells.Item(2, 13).Value = -196.33334
$ws.Cells.Item(2, 14).Value = -612.6
$ws.Cells.Item(97, 8).Value = 744
$ws.Cells.Item(97, 9).Value = 725.3333
$ws.Cells.Item(97, 11).Value = 725.3333
$ws.Cells.Item(97, 13).Value = -229.3333
$ws.Cells.Item(102, 8).Value = 1440.65
$ws.Cells.Item(102, 9).Value = 1289.6666
$ws.Cells.Item(102, 11).Value = 1289.6666
$ws.Cells.Item(102, 13).Value = 332.3334

# --- LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(2, 8).Value = 29750
$ws.Cells.Item(2, 9).Value = 3666.6667
$ws.Cells.Item(2, 10).Value = 108000
$ws.Cells.Item(2, 11).Value = 3666.6667
$ws.Cells.Item(2, 12).Value = 108000
$ws.Cells.Item(2, 13).Value = -3554.6667
$ws.Cells.Item(2, 14).Value = -108224
$ws.Cells.Item(22, 8).Value = 680.1667
$ws.Cells.Item(22, 9).Value = 658.4
$ws.Cells.Item(22, 11).Value = 658.4
$ws.Cells.Item(22, 13).Value = -363.4
$ws.Cells.Item(24, 8).Value = 1001333.3
$ws.Cells.Item(24, 9).Value = 1001333.3
$ws.Cells.Item(24, 11).Value = 1001333.3
$ws.Cells.Item(24, 13).Value = -1000990.3
$ws.Cells.Item(27, 8).Value = 680.1667
$ws.Cells.Item(27, 9).Value = 658.4
$ws.Cells.Item(27, 11).Value = 658.4
$ws.Cells.Item(27, 13).Value = -551.4
$ws.Cells.Item(55, 8).Value = 1203.375
$ws.Cells.Item(55, 9).Value = 1203.375
$ws.Cells.Item(55, 11).Value = 1203.375
$ws.Cells.Item(55, 13).Value = -1030.375

# --- WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(132, 8).Value = 1783.7894
$ws.Cells.Item(132, 9).Value = 1716.2222
$ws.Cells.Item(132, 11).Value = 5148.6666
$ws.Cells.Item(132, 13).Value = -2618.6666
